$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.956.02'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '2.213.86'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('E4').Value = '  -0.14%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '256.13'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.31%  '
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '77.57'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +3.44%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.612'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('E8').Value = '  -0.09%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.598'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.24%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '42.83'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +3.24%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0912'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.66%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '7.04'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +2.05%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.103'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.27%  '
$ws.Range('D14').Value = '2.544.78'
$ws.Range('E14').Value = '  -1.23%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '14.42'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').Value = '2.210.72'
$ws.Range('E16').Value = '  -1.19%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.781'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').Value = '42.909.93'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('E19').Value = '  -0.25%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '71.15'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.03%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '5.98'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.07%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '2.33'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +7.20%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '230.20'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.27%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '9.28'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -3.97%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '42.67'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +8.45%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '10.77'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.34%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '3.35'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -3.06%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.20'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -2.57%  '
$ws.Range('E30').Value = '  -1.06%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '172.69'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '20.39'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.73%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.0876'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +9.52%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '5.21'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.80%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.122'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.05%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.0359'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +8.82%  '
$ws.Range('E37').Value = '  -2.41%  '
$ws.Range('E38').Value = '  -3.04%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '13.22'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +2.63%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.94'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +21.00%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.11'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.48%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.202'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -1.41%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '61.34'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.99%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '5.31'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.34%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '103.46'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.08%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '8.48'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -2.24%  '
$ws.Range('B47').Value = 'WOONetwork'
$ws.Range('C47').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.469'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0975'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -1.43%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.12'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.14'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.50'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +24.32%  '
